# "fixed dip switch tables"
#
# The DIP Switches worksheet has three small reference tables describing
# what each DIP switch position means. This fixes incorrect / duplicated
# text in those tables:
#   - Table 1 (rows 2-11): switch 1 row was showing "Minus"/" Plus" (left
#     over from the Counter-mode table) but should read "not used" for
#     the Timer-mode table.
#   - Table 3 (rows 24-33): the hour-modification rows used inconsistent
#     "On"/"Off" casing, and the 2 Hour / 1 Hour rows had the OFF column
#     mistakenly duplicating the ON text instead of showing the real OFF
#     text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DIP Switches")

# --- Table 1 (rows 2-11): switch 1 is not used for this table ---
$ws.Range("C4").Value = "not used"
$ws.Range("D4").Value = "not used"

# --- Table 3 (rows 24-33): hour modification ON/OFF fixes ---
$ws.Range("C27").Value = "8 Hour Modification ON"
$ws.Range("D27").Value = "8 Hour Modification OFF"

$ws.Range("C28").Value = "4 Hour Modification ON"
$ws.Range("D28").Value = "4 Hour Modification OFF"

$ws.Range("C29").Value = "2 Hour Modification ON"
$ws.Range("D29").Value = "2 Hour Modification OFF"

$ws.Range("C30").Value = "1 Hour Modification ON"
$ws.Range("D30").Value = "1 Hour Modification OFF"

# D2/D13 had a stray "applyAlignment" style with no effect; clear it back
# to the Normal style so it matches the rest of the header row.
$ws.Range("D2").Style = "Normal"
$ws.Range("D13").Style = "Normal"

# Restore the sheet's active selection to D40 (matches the author's saved
# view state).
$ws.Activate()
$ws.Range("D40").Select()
